$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 144.4
$ws.Range("C3").Value = 163.2
$ws.Range("C4").Value = 162.4
$ws.Range("C6").Value = 184
$ws.Range("C8").Value = 171.2
$ws.Range("C18").Value = 166.3
